# DT added label option
# Re-ran the backward-elimination regression export on a later date/time;
# refresh the "Date:" / "Time:" stamps embedded in each worksheet's
# statsmodels OLS summary text (cell B2 on every sheet) to match the new run.

$wb = $excel.ActiveWorkbook

$oldDate = "Sun, 05 Jan 2020"
$newDate = "Wed, 08 Jan 2020"
$oldTimes = @("21:22:09", "21:22:10")

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $newTime = "19:07:15"
    if ($i -gt 20) {
        $newTime = "19:07:16"
    }

    $cell = $ws.Range("B2")
    $text = $cell.Value()

    if ($text -ne $null) {
        $text = $text.Replace($oldDate, $newDate)
        foreach ($ot in $oldTimes) {
            $text = $text.Replace($ot, $newTime)
        }
        $cell.Value = $text
    }
}

Write-Output "done"
